$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Fitness column (C) for rows 2 through 252 is updated to a constant
# value of 7310 (previously it held several different plateau values).
$ws.Range("C2:C252").Value = 7310
